# Made fixes for Prod/Demo Verification Script and Display Convenience fees qa done for all versions
$wb = $excel.ActiveWorkbook

$wsProfile = $wb.Worksheets.Item("CreateModifyDeleteProfile")
$wsProfile.Range("B2").Value = "Wed May 21 14:42:34 IST 2025"
$wsProfile.Range("C2").Value = "Pass"

$wsCC = $wb.Worksheets.Item("AddModifyDeleteCC")
$wsCC.Range("B2").Value = "Wed May 21 15:57:03 IST 2025"
$wsCC.Range("C2").Value = "Pass"

$wsACH = $wb.Worksheets.Item("AddModifyDeleteACH")
$wsACH.Range("B2").Value = "Wed May 21 15:52:36 IST 2025"
$wsACH.Range("B3").Value = "Wed May 21 15:53:50 IST 2025"
$wsACH.Range("B4").Value = "Wed May 21 15:54:57 IST 2025"
